$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 12412.92
$ws.Range("I138").Value = 7329.3335
$ws.Range("J138").Value = 12737.404
$ws.Range("K138").Value = 21988.0005
$ws.Range("L138").Value = 38212.212
$ws.Range("M138").Value = -16848.0005
$ws.Range("N138").Value = -48492.212

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1693.5172
$ws.Range("I2").Value = 1171.7
$ws.Range("J2").Value = 2853.111
$ws.Range("K2").Value = 1171.7
$ws.Range("L2").Value = 2853.111
$ws.Range("M2").Value = -1058.7
$ws.Range("N2").Value = -3079.111
$ws.Range("H32").Value = 13702221
$ws.Range("I32").Value = 15153454
$ws.Range("K32").Value = 15153454
$ws.Range("M32").Value = -15153167
$ws.Range("H45").Value = 898
$ws.Range("I45").Value = 898
$ws.Range("K45").Value = 898
$ws.Range("M45").Value = -521
$ws.Range("H61").Value = 66670130
$ws.Range("I61").Value = 71431990
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 71431990
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -71431778
$ws.Range("N61").Value = -4424
$ws.Range("H97").Value = 1933.6316
$ws.Range("I97").Value = 2016.2222
$ws.Range("J97").Value = 447
$ws.Range("K97").Value = 2016.2222
$ws.Range("L97").Value = 447
$ws.Range("M97").Value = -1520.2222
$ws.Range("N97").Value = -1439
$ws.Range("H110").Value = 19889.191
$ws.Range("I110").Value = 21827.525
$ws.Range("J110").Value = 1475
$ws.Range("K110").Value = 21827.525
$ws.Range("L110").Value = 1475
$ws.Range("M110").Value = -19782.525
$ws.Range("N110").Value = -5565
$ws.Range("H116").Value = 1693.5172
$ws.Range("I116").Value = 1171.7
$ws.Range("J116").Value = 2853.111
$ws.Range("K116").Value = 1171.7
$ws.Range("L116").Value = 2853.111
$ws.Range("M116").Value = 1122.3
$ws.Range("N116").Value = -7441.111
$ws.Range("H122").Value = 7409650
$ws.Range("I122").Value = 1836.2285
$ws.Range("J122").Value = 33337000
$ws.Range("K122").Value = 5508.6855
$ws.Range("L122").Value = 100011000
$ws.Range("M122").Value = -3058.6855
$ws.Range("N122").Value = -100015900
$ws.Range("H132").Value = 52711444
$ws.Range("I132").Value = 10380.556
$ws.Range("J132").Value = 100142400
$ws.Range("K132").Value = 31141.668
$ws.Range("L132").Value = 300427200
$ws.Range("M132").Value = -28611.668
$ws.Range("N132").Value = -300432260
$ws.Range("H136").Value = 66670130
$ws.Range("I136").Value = 71431990
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 214295970
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -214293420
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1693.5172
$ws.Range("I3").Value = 1171.7
$ws.Range("J3").Value = 2853.111
$ws.Range("K3").Value = 1171.7
$ws.Range("L3").Value = 2853.111
$ws.Range("M3").Value = -1057.7
$ws.Range("N3").Value = -3081.111
$ws.Range("H94").Value = 1620.3636
$ws.Range("I94").Value = 928.9474
$ws.Range("J94").Value = 5999.3335
$ws.Range("K94").Value = 928.9474
$ws.Range("L94").Value = 5999.3335
$ws.Range("M94").Value = -477.9474
$ws.Range("N94").Value = -6901.3335
$ws.Range("H105").Value = 12520.077
$ws.Range("I105").Value = 31080.25
$ws.Range("J105").Value = 4271.1113
$ws.Range("K105").Value = 31080.25
$ws.Range("L105").Value = 4271.1113
$ws.Range("M105").Value = -29333.25
$ws.Range("N105").Value = -7765.1113
$ws.Range("H107").Value = 6806.6
$ws.Range("I107").Value = 2020.5
$ws.Range("K107").Value = 2020.5
$ws.Range("M107").Value = -100.5
$ws.Range("H134").Value = 2434.9487
$ws.Range("I134").Value = 2269.2974
$ws.Range("K134").Value = 6807.8922
$ws.Range("M134").Value = -4272.8922

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32054792
$ws.Range("I31").Value = 3276.6
$ws.Range("J31").Value = 43107040
$ws.Range("K31").Value = 3276.6
$ws.Range("L31").Value = 43107040
$ws.Range("M31").Value = -2981.6
$ws.Range("N31").Value = -43107630
$ws.Range("H34").Value = 32054792
$ws.Range("I34").Value = 3276.6
$ws.Range("J34").Value = 43107040
$ws.Range("K34").Value = 3276.6
$ws.Range("L34").Value = 43107040
$ws.Range("M34").Value = -3074.6
$ws.Range("N34").Value = -43107444

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27749118
$ws.Range("I4").Value = 48856464
$ws.Range("J4").Value = 13470619
$ws.Range("K4").Value = 146569392
$ws.Range("L4").Value = 40411857
$ws.Range("M4").Value = -146569280
$ws.Range("N4").Value = -40412081
$ws.Range("H40").Value = 631.5
$ws.Range("I40").Value = 46.666668
$ws.Range("J40").Value = 882.1429000000001
$ws.Range("K40").Value = 186.666672
$ws.Range("L40").Value = 3528.5716
$ws.Range("M40").Value = -117.666672
$ws.Range("N40").Value = -3666.5716
$ws.Range("H132").Value = 4765820
$ws.Range("I132").Value = 1519.8
$ws.Range("K132").Value = 13678.2
$ws.Range("M132").Value = -11148.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 977.2941
$ws.Range("I97").Value = 927.6667
$ws.Range("K97").Value = 927.6667
$ws.Range("M97").Value = -431.6667
$ws.Range("H122").Value = 27780746
$ws.Range("I122").Value = 2228.9333
$ws.Range("J122").Value = 166673330
$ws.Range("K122").Value = 6686.7999
$ws.Range("L122").Value = 500019990
$ws.Range("M122").Value = -4236.7999
$ws.Range("N122").Value = -500024890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1982.2424
$ws.Range("I46").Value = 1018.2917
$ws.Range("J46").Value = 4552.778
$ws.Range("K46").Value = 1018.2917
$ws.Range("L46").Value = 4552.778
$ws.Range("M46").Value = -830.2917
$ws.Range("N46").Value = -4928.778
$ws.Range("H61").Value = 3337.5715
$ws.Range("I61").Value = 2378.8635
$ws.Range("J61").Value = 6852.8335
$ws.Range("K61").Value = 2378.8635
$ws.Range("L61").Value = 6852.8335
$ws.Range("M61").Value = -2176.8635
$ws.Range("N61").Value = -7256.8335
$ws.Range("H93").Value = 693134.25
$ws.Range("I93").Value = 3813.2778
$ws.Range("J93").Value = 3795078.5
$ws.Range("K93").Value = 3813.2778
$ws.Range("L93").Value = 3795078.5
$ws.Range("M93").Value = -2565.2778
$ws.Range("N93").Value = -3797574.5
$ws.Range("H113").Value = 3337.5715
$ws.Range("I113").Value = 2378.8635
$ws.Range("J113").Value = 6852.8335
$ws.Range("K113").Value = 2378.8635
$ws.Range("L113").Value = 6852.8335
$ws.Range("M113").Value = -208.8634999999999
$ws.Range("N113").Value = -11192.8335
$ws.Range("H122").Value = 3208871.2
$ws.Range("I122").Value = 3418.742
$ws.Range("K122").Value = 10256.226
$ws.Range("M122").Value = -7806.226000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5364.4287
$ws.Range("I96").Value = 4758.5
$ws.Range("J96").Value = 9000
$ws.Range("K96").Value = 4758.5
$ws.Range("L96").Value = 9000
$ws.Range("M96").Value = -3385.5
$ws.Range("N96").Value = -11746
$ws.Range("H113").Value = 976.15
$ws.Range("I113").Value = 959.3570999999999
$ws.Range("J113").Value = 1015.3333
$ws.Range("K113").Value = 2878.0713
$ws.Range("L113").Value = 3045.9999
$ws.Range("M113").Value = -708.0712999999996
$ws.Range("N113").Value = -7385.9999

Write-Output "Applied 189 cell updates across 8 sheets"